$d = $word.ActiveDocument

# --- 1. Rename CART3_FAC.slx -> CART3_FIRC.slx everywhere it appears verbatim ---
# (covers: Title paragraph, Introduction italic mention, "Running CART3_FAC.slx" heading)
$rng = $d.Content
$rng.Find.Execute("CART3_FAC.slx", $true, $false, $false, $false, $false, $true, 1, $false, "CART3_FIRC.slx", 2)

# --- 2. Rename CART3_FAC_paper -> CART3_FIRC_paper (both occurrences) ---
$rng = $d.Content
$rng.Find.Execute("CART3_FAC_paper", $true, $false, $false, $false, $false, $true, 1, $false, "CART3_FIRC_paper", 2)

# --- 3. "add FAC functionality" -> "add FIRC functionality" ---
$rng = $d.Content
$rng.Find.Execute("add FAC functionality", $true, $false, $false, $false, $false, $true, 1, $false, "add FIRC functionality", 2)

# --- 4. Replace the paper-title sentence in the Introduction paragraph ---
# Old (spans 3 runs, only middle italic run + trailing "." + italic space):
#   [italic] Fault Adaptive Wind Turbine Controller | [plain] . | [italic] (space)
# New (all italic):
#   Demonstration of a Fault Impact Reduction Control Module for Wind Turbines. (space)
$rng = $d.Content
$rng.Find.Execute("Fault Adaptive Wind Turbine Controller", $true, $false, $false, $false, $false, $true, 1, $false, "Demonstration of a Fault Impact Reduction Control Module for Wind Turbines.", 2)

# Now remove the old trailing non-italic "." run and make sure the italic run picks up
# the rest (the italic space run following it stays, so we just delete the stray ".").
$rng = $d.Content
$found = $rng.Find.Execute("Wind Turbines.")
if ($found) {
    $period = $d.Range($rng.End, $rng.End + 1)
    if ($period.Text -eq ".") {
        $period.Italic = 1
        $period.Delete()
    }
}

Write-Host "done"
